# Entrevista.docx — font-size update (commit: "Atualizacao de Tamanho de Fonte /
# Mudanca da Fonte para 14") plus the accompanying _GoBack bookmark
# repositioning and stale page-break-cache cleanup that Word performs when
# it resaves a document after an in-place edit near that location.

$d = $word.ActiveDocument

# --- 1. Title paragraph ("Entrevista 1") : 72pt -> 20pt (sz/szCs 144 -> 40)
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Size = 20
$p1.Range.Font.SizeBi = 20

# --- 2. First blank line under the title : 12pt -> 14pt (sz/szCs 24 -> 28)
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Font.Size = 14
$p2.Range.Font.SizeBi = 14

# Paragraph 3 (the next blank line) is left untouched (stays 12pt / sz 24).

# --- 3. "-Perguntar se ele quer implementar as comandas..." paragraph : 72pt -> 14pt
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Font.Size = 14
$p4.Range.Font.SizeBi = 14

# --- 4. Blank line afterwards : 72pt -> 14pt
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Font.Size = 14
$p5.Range.Font.SizeBi = 14

# --- 5. "-Quais tipos de " paragraph : 72pt -> 14pt
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Font.Size = 14
$p6.Range.Font.SizeBi = 14

# Re-run the text through Find/Replace so Word drops the stale
# <w:lastRenderedPageBreak/> rendering-cache marker that was sitting in
# front of this run (Word recomputes it on layout and strips it here).
$d.Content.Find.Execute("Quais tipos de", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Quais tipos de", 2) | Out-Null

# --- 6. Move the "_GoBack" bookmark so it again wraps the very last edit,
# i.e. the title paragraph ("Entrevista 1"), the way Word stamps it after
# a save that touched that text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$titleRange = $d.Paragraphs.Item(1).Range
$d.Bookmarks.Add("_GoBack", $titleRange) | Out-Null

Write-Output "done"
